# Insert a new worksheet "Hospital Santa Rita" right after
# "Hospital Pequeno Príncipe" (currently sheet #11) and populate it with
# the program/vacancy/applicant/competition-ratio table.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item(11)
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Hospital Santa Rita"

$data = @(
    @("PROGRAMA", "VAGAS", "INSCRITOS", "CONCORRÊNCIA"),
    @("Anestesiologia", "3", "80", "26,67"),
    @("Cirurgia Geral", "4", "98", "24,5"),
    @("Clínica Médica", "4", "74", "18,5"),
    @("Infectologia", "1", "10", "10"),
    @("Medicina Intensivo Adulto", "2", "16", "8"),
    @("Ortopedia e Traumatologia", "3", "49", "16,33"),
    @("Pediatria", "2", "44", "22"),
    @("Radiologia e Diagnóstico por Imagem", "2", "47", "23,5"),
    @("Cardiologia", "2", "6", "3"),
    @("Cirurgia Oncológica", "1", "6", "6"),
    @("Cirurgia Vascular", "1", "6", "6"),
    @("Coloproctologia", "1", "6", "6"),
    @("Ecocardiografia", "1", "4", "4"),
    @("Gastroenterologia", "1", "4", "4"),
    @("Neonatologia", "2", "4", "2"),
    @("Oncologia Clínica", "1", "3", "3")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Match the header-row formatting used by every other sheet in the
# workbook: bold, centered horizontally, top-aligned vertically, with a
# thin box border.
$header = $ws.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
for ($i = 1; $i -le 4; $i++) {
    $header.Borders.Item($i).LineStyle = 1
    $header.Borders.Item($i).Weight = 2
}

$ws.Range("A1").Select()
